$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85, pushing existing rows 85..125 down to 86..126
$ws.Rows.Item(85).Insert()

# Fill in the new row 85 - same as what used to be row 85 (now row 86),
# except for the date and price columns which carry new figures.
$ws.Cells.Item(85, 1).Value = 11
$ws.Cells.Item(85, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(85, 3).Value = "Bíobío"
$ws.Cells.Item(85, 4).Value = 44680
$ws.Cells.Item(85, 5).Value = 8
$ws.Cells.Item(85, 6).Value = 100112043
$ws.Cells.Item(85, 7).Value = "Pepino ensalada"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 100
$ws.Cells.Item(85, 11).Value = 15000
$ws.Cells.Item(85, 12).Value = 16000
$ws.Cells.Item(85, 13).Value = 15500
$ws.Cells.Item(85, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 258
$ws.Cells.Item(85, 17).Value = 60
$ws.Cells.Item(85, 18).Value = "Hortaliza"
